$d = $word.ActiveDocument
$full = $d.Content
$frag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>-Creator:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Windows Applikation, welche die Messdaten </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">des Device-Locator empfängt und verarbeitet. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Die Applikation </w:t>
      </w:r>
      <w:r>
        <w:t>öffnet ein Fenster mit einem GUI und einer Visualisierung.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Links sind </w:t>
      </w:r>
      <w:r>
        <w:t>Buttons,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> mit welchen man verschiedene Parameter verändern kann oder bestimmte Funktionen ausführen kann. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Rechts ist eine Visualisierung der Messpunkte, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmaps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> und dem Raumlayout.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Ein blau blinkendes Rechteck zeigt die Position des aktuellen Messpunktes an. Man kann die Position per Mausklick ändern oder über die Steuersignale die der esp32 sendet.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Links oben werden 2 Werte angezeigt. Der </w:t>
      </w:r>
      <w:r>
        <w:t>linke Wert</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> zeigt die Signalstärke des ausgewählten Messpunktes</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, der aktuellen </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> an</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Der rechte Wert gibt die „Qualität“ der aktuellen </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> an. Mit dem ersten Button kann man zwischen der Datenpunkt- und der </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmapansicht</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> wechseln. </w:t>
      </w:r>
      <w:r>
        <w:t>Der nächste Button „</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> gen.“ berechnet alle </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmaps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> aus den aktuellen Messpunkten</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Der Button „</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Loeschen</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">“ entfernt alle Messpunkte. </w:t>
      </w:r>
      <w:r>
        <w:t>Der Button</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> „Speichern“ speichert die aktuellen Messpunkte in einer Datei mit dem aktuellen Zeitstempel. </w:t>
      </w:r>
      <w:r>
        <w:t>Der Button „Laden“</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>lädt eine Datei mit dem Namen „</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">“ und generiert alle </w:t>
      </w:r>
      <w:r>
        <w:t>Messpunkte,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> die in dieser gespeichert waren. </w:t>
      </w:r>
      <w:r>
        <w:t>Der nächste Button wechselt die aktuellen Messpunktindexe durch, welche durch die Anzeigefunktionen gebraucht werden.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Der nächste Button legt fest, ob das Programm einen empfangenen Messpunkt als </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">neuen an der aktuellen Position hinzufügen soll oder ob diese zur Suche benutzt werden. Der Wert und die 2 </w:t>
      </w:r>
      <w:r>
        <w:t>Knöpfe,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> um diesen zu erhöhen/verringern ist dazu da, um die generierte </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Suchheatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> anzupassen. Die </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Suchheatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> visualisiert die Wahrscheinlichkeit, dass der esp32 sich an der entsprechenden Position befindet, mit höheren Werten, kann man geringere Wahrscheinlichkeiten weniger intensiv darstellen, womit man besser erkennen kann, wo das Gerät sich wahrscheinlich befindet (es dient nur für die bessere Visualisierung). Der nächste Button legt fest, ob nur die aktuelle </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmap</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> zur Suche benutzt wird oder alle </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmaps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. Der Button „Gewichtung an/aus“ legt fest, ob qualitativ hochwertigere </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmaps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> mehr in die Positionsbestimmung mit einberechnet, werden sollen oder nicht. Diese Option kann zu besseren Ergebnissen führen, falls es qualitativ schlechte </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Heatmaps</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> gibt. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Device-Locator:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Ein Esp32 Programm, welches Messdaten</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> und Steuersignale </w:t>
      </w:r>
      <w:r>
        <w:t>verschicket und</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> per Netzwerk konfiguriert werden kann.</w:t>
      </w:r>
    </w:p>
    </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void]$full.InsertXML($frag)
